# Generate Report for Handback
# Fills in the "b7629e25-9312-4863-ba80-81cf67b25162" row (row 7) on the
# zh-cn and de-de sheets with the result of a handback attempt that failed
# because the handback file version was stale.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/615f4b2694c78570e81a1c173793726e35109a47/e2e/b7629e25-9312-4863-ba80-81cf67b25162.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/272e74c0d1bc82a08cd5513be66961488b876c71/e2e/b7629e25-9312-4863-ba80-81cf67b25162.md."
$latestMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/272e74c0d1bc82a08cd5513be66961488b876c71/e2e/b7629e25-9312-4863-ba80-81cf67b25162.md"
$mdDisplay = "b7629e25-9312-4863-ba80-81cf67b25162.md"

# zh-cn sheet (row 7 : b7629e25-9312-4863-ba80-81cf67b25162)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("I7").Value = $mdDisplay
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $latestMdUrl, "", "", $mdDisplay)
$wsZh.Range("I7").Style = "HyperLink"
$wsZh.Range("J7").Value = "b7629e25-9312-4863-ba80-81cf67b25162.af35b29bba625a195395e37141acdadae7f22364.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-16 14:56:18"
$wsZh.Range("P7").Value = $errorDetail

# de-de sheet (row 7 : b7629e25-9312-4863-ba80-81cf67b25162)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("I7").Value = $mdDisplay
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $latestMdUrl, "", "", $mdDisplay)
$wsDe.Range("I7").Style = "HyperLink"
$wsDe.Range("J7").Value = "b7629e25-9312-4863-ba80-81cf67b25162.af35b29bba625a195395e37141acdadae7f22364.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-16 14:56:25"
$wsDe.Range("P7").Value = $errorDetail
